$wb = $excel.ActiveWorkbook

# --- LogInConfiguration: move selection from B15 to D2 ---
$wsLogin = $wb.Worksheets.Item("LogInConfiguration")
$wsLogin.Range("D2").Select()

# --- ReportConfiguration: move selection from C22 to C2 (also scrolled topLeftCell C1 -> B1) ---
$wsReport = $wb.Worksheets.Item("ReportConfiguration")
$wsReport.Range("C2").Select()

# --- Welcome: selection moves from single cell C9 to range C2:C9 (active cell C2); tab no longer active ---
$wsWelcome = $wb.Worksheets.Item("Welcome")
$wsWelcome.Range("C2:C9").Select()

# --- FindFormRT: update CanExecute values from "n" to "y" for rows 11-16 and 18 (not row 17) ---
$wsFind = $wb.Worksheets.Item("FindFormRT")
$wsFind.Range("C11").Value = "y"
$wsFind.Range("C12").Value = "y"
$wsFind.Range("C13").Value = "y"
$wsFind.Range("C14").Value = "y"
$wsFind.Range("C15").Value = "y"
$wsFind.Range("C16").Value = "y"
$wsFind.Range("C18").Value = "y"

# --- FindFormRT becomes the active sheet/tab; selection (within frozen pane) moves from E10 to C19 ---
$wsFind.Activate()
$wsFind.Range("C19").Select()
